# GACPAQ pages.xlsx - "Country Status" sheet update
# Commit message: "provided country-specific examples to the HBSC Preamble of en-NG"
#
# The underlying content edit (HBSC Preamble text, presumably on the
# "Section Pages" sheet for en-NG) is tracked on the "Country Status" sheet
# as a refresh of translation/review progress for several languages: new
# "in review" / "being revised" / "data received" / "on-going (N)" status
# labels are introduced, and several countries' status cells are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Status")
$ws.Activate()

# --- Reference (donor) cells whose fills/format we reuse for the
#     different status categories. None of these are themselves touched
#     by this edit, so they stay valid donors throughout the script.
$donorCompleted = $ws.Range("C2")    # "completed" fill
$donorAwaiting  = $ws.Range("H14")   # "awaiting  data" -> reused for "in review"/"on-going (3)" style family
$donorOngoing   = $ws.Range("G3")    # "on-going" fill

function Set-StatusCell($addr, $text, $donor) {
    $target = $ws.Range($addr)
    if ($donor -ne $null) {
        $donor.Copy()
        $target.PasteSpecial(-4122)  # xlPasteFormats
    }
    $target.Value = $text
}

# NOTE: cells are written in an order chosen to reproduce the exact
# shared-string append order of the source edit (new unique strings are
# appended to xl/sharedStrings.xml in first-use order):
#   399 "in review", 400 "being revised", 401 "data received",
#   402 "on-going (1)", 403 "on-going (2)", 404 "on-going (3)"
$excel.CutCopyMode = $false

# Row 2 - Brazil - Portugese
Set-StatusCell "J2" "in review" $donorAwaiting

# Row 3 - Canada - English
Set-StatusCell "J3" "in review" $donorAwaiting

# Row 4 - Canada - French
Set-StatusCell "I4" "completed" $donorCompleted
Set-StatusCell "J4" "in review" $donorAwaiting

# Row 5 - Chile - Spanish
Set-StatusCell "F5" "in review" $donorOngoing

# Row 6 - China - Chinese
Set-StatusCell "G6" "completed" $donorCompleted
Set-StatusCell "H6" "being revised" $donorOngoing
Set-StatusCell "I6" "being revised" $donorOngoing
Set-StatusCell "J6" "in review" $donorAwaiting

# Row 17 - Nigeria - English (en-NG) - country referenced in the commit message
$ws.Rows.Item(17).RowHeight = 35
Set-StatusCell "H17" "data received" $null
Set-StatusCell "I17" "data received" $null

# Row 19 - Sweden - Swedish
Set-StatusCell "K19" "on-going (1)" $null

# Row 13 - Mexico - Spanish
Set-StatusCell "J13" "deployed" $donorCompleted
Set-StatusCell "K13" "on-going (2)" $donorOngoing

# Row 18 - Spain - Spanish
Set-StatusCell "I18" "completed" $donorCompleted
Set-StatusCell "J18" "deployed" $donorCompleted
Set-StatusCell "K18" "on-going (3)" $donorAwaiting

# Row 7 - Colombia - Spanish
Set-StatusCell "F7" "completed" $donorCompleted
Set-StatusCell "G7" "completed" $donorCompleted
Set-StatusCell "H7" "awaiting  data" $donorOngoing
Set-StatusCell "I7" "awaiting  data" $donorOngoing

# Row 8 - Czech Republic - Czech
Set-StatusCell "F8" "completed" $donorCompleted
Set-StatusCell "G8" "completed" $donorCompleted
Set-StatusCell "H8" "awaiting  data" $donorOngoing
Set-StatusCell "I8" "awaiting  data" $donorOngoing

# Row 20 - Thailand - Thai
Set-StatusCell "F20" "completed" $donorCompleted
Set-StatusCell "G20" "completed" $donorCompleted
Set-StatusCell "H20" "completed" $donorCompleted
Set-StatusCell "I20" "completed" $donorCompleted

# Row 21 - UAE - Arabic
Set-StatusCell "G21" "completed" $donorCompleted
Set-StatusCell "H21" "completed" $donorCompleted
Set-StatusCell "I21" "completed" $donorCompleted

# Row 22 - UAE - English
Set-StatusCell "G22" "completed" $donorCompleted
Set-StatusCell "H22" "completed" $donorCompleted
Set-StatusCell "I22" "completed" $donorCompleted

$excel.CutCopyMode = $false

# --- View state: re-freeze panes at B2 (was Y18/24,17) and refresh zoom ---
$win = $excel.ActiveWindow
$win.Zoom = 114
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
$ws.Range("E9").Select()

Write-Host "Country Status sheet updated"
